$wb = $excel.ActiveWorkbook

# --- 1) Update "总计" (Sheet1): insert new 2022-Q4 row at top, shift others down ---
$total = $wb.Worksheets.Item(1)
for ($r = 6; $r -ge 2; $r--) {
    $srcRow = $total.Range("A" + $r + ":D" + $r)
    $dstRow = $total.Range("A" + ($r+1) + ":D" + ($r+1))
    $srcRow.Copy()
    $dstRow.PasteSpecial(-4163)
}
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 27
$total.Range("D2").Value = 0.73
for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
# newly-extended row 7 does not inherit the style copied for A2:A6 via PasteSpecial; fix it up
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

# --- 2) Insert new worksheet "2022-Q4" right after "总计", before "2021-Q4" ---
$refSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($refSheet)
$newSheet.Name = "2022-Q4"

# Copy the header/data cell styling pattern (s="2" on header row + column A) from an existing
# quarter sheet that already has enough rows, so the new sheet matches the workbook's look.
$styleSrc = $wb.Worksheets.Item("2021-Q2")
$styleSrc.Range("A1:H28").Copy()
$newSheet.Range("A1:H28").PasteSpecial(-4122)

# --- 3) Fill header row ---
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# --- 4) Fill data rows ---
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "012988"
$newSheet.Cells.Item(2, 3).Value = "嘉合锦明混合C"
$newSheet.Cells.Item(2, 4).Value = "3.26"
$newSheet.Cells.Item(2, 5).Value = "87.23"
$newSheet.Cells.Item(2, 6).Value = "3.81"
$newSheet.Cells.Item(2, 7).Value = "0.1242"
$newSheet.Cells.Item(2, 8).Value = 8
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "012987"
$newSheet.Cells.Item(3, 3).Value = "嘉合锦明混合A"
$newSheet.Cells.Item(3, 4).Value = "2.08"
$newSheet.Cells.Item(3, 5).Value = "87.23"
$newSheet.Cells.Item(3, 6).Value = "3.81"
$newSheet.Cells.Item(3, 7).Value = "0.0792"
$newSheet.Cells.Item(3, 8).Value = 8
$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "011977"
$newSheet.Cells.Item(4, 3).Value = "格林研究优选混合A"
$newSheet.Cells.Item(4, 4).Value = "2.10"
$newSheet.Cells.Item(4, 5).Value = "86.76"
$newSheet.Cells.Item(4, 6).Value = "3.64"
$newSheet.Cells.Item(4, 7).Value = "0.0764"
$newSheet.Cells.Item(4, 8).Value = 10
$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(5, 2).Value = "006973"
$newSheet.Cells.Item(5, 3).Value = "太平睿盈混合A"
$newSheet.Cells.Item(5, 4).Value = "3.17"
$newSheet.Cells.Item(5, 5).Value = "29.39"
$newSheet.Cells.Item(5, 6).Value = "1.82"
$newSheet.Cells.Item(5, 7).Value = "0.0577"
$newSheet.Cells.Item(5, 8).Value = 5
$newSheet.Cells.Item(6, 1).Value = 4
$newSheet.Cells.Item(6, 2).Value = "006424"
$newSheet.Cells.Item(6, 3).Value = "嘉合锦程价值精选混合A"
$newSheet.Cells.Item(6, 4).Value = "1.40"
$newSheet.Cells.Item(6, 5).Value = "82.04"
$newSheet.Cells.Item(6, 6).Value = "3.83"
$newSheet.Cells.Item(6, 7).Value = "0.0536"
$newSheet.Cells.Item(6, 8).Value = 10
$newSheet.Cells.Item(7, 1).Value = 5
$newSheet.Cells.Item(7, 2).Value = "700001"
$newSheet.Cells.Item(7, 3).Value = "平安行业先锋混合"
$newSheet.Cells.Item(7, 4).Value = "1.77"
$newSheet.Cells.Item(7, 5).Value = "92.09"
$newSheet.Cells.Item(7, 6).Value = "2.99"
$newSheet.Cells.Item(7, 7).Value = "0.0529"
$newSheet.Cells.Item(7, 8).Value = 10
$newSheet.Cells.Item(8, 1).Value = 6
$newSheet.Cells.Item(8, 2).Value = "015011"
$newSheet.Cells.Item(8, 3).Value = "嘉合锦鑫混合C"
$newSheet.Cells.Item(8, 4).Value = "1.53"
$newSheet.Cells.Item(8, 5).Value = "83.79"
$newSheet.Cells.Item(8, 6).Value = "3.26"
$newSheet.Cells.Item(8, 7).Value = "0.0499"
$newSheet.Cells.Item(8, 8).Value = 9
$newSheet.Cells.Item(9, 1).Value = 7
$newSheet.Cells.Item(9, 2).Value = "008261"
$newSheet.Cells.Item(9, 3).Value = "招商研究优选股票A"
$newSheet.Cells.Item(9, 4).Value = "1.90"
$newSheet.Cells.Item(9, 5).Value = "86.08"
$newSheet.Cells.Item(9, 6).Value = "2.57"
$newSheet.Cells.Item(9, 7).Value = "0.0488"
$newSheet.Cells.Item(9, 8).Value = 7
$newSheet.Cells.Item(10, 1).Value = 8
$newSheet.Cells.Item(10, 2).Value = "015010"
$newSheet.Cells.Item(10, 3).Value = "嘉合锦鑫混合A"
$newSheet.Cells.Item(10, 4).Value = "1.37"
$newSheet.Cells.Item(10, 5).Value = "83.79"
$newSheet.Cells.Item(10, 6).Value = "3.26"
$newSheet.Cells.Item(10, 7).Value = "0.0447"
$newSheet.Cells.Item(10, 8).Value = 9
$newSheet.Cells.Item(11, 1).Value = 9
$newSheet.Cells.Item(11, 2).Value = "006425"
$newSheet.Cells.Item(11, 3).Value = "嘉合锦程价值精选混合C"
$newSheet.Cells.Item(11, 4).Value = "1.02"
$newSheet.Cells.Item(11, 5).Value = "82.04"
$newSheet.Cells.Item(11, 6).Value = "3.83"
$newSheet.Cells.Item(11, 7).Value = "0.0391"
$newSheet.Cells.Item(11, 8).Value = 10
$newSheet.Cells.Item(12, 1).Value = 10
$newSheet.Cells.Item(12, 2).Value = "007669"
$newSheet.Cells.Item(12, 3).Value = "太平睿盈混合C"
$newSheet.Cells.Item(12, 4).Value = "0.99"
$newSheet.Cells.Item(12, 5).Value = "29.39"
$newSheet.Cells.Item(12, 6).Value = "1.82"
$newSheet.Cells.Item(12, 7).Value = "0.0180"
$newSheet.Cells.Item(12, 8).Value = 5
$newSheet.Cells.Item(13, 1).Value = 11
$newSheet.Cells.Item(13, 2).Value = "005493"
$newSheet.Cells.Item(13, 3).Value = "鑫元价值精选灵活配置混合A"
$newSheet.Cells.Item(13, 4).Value = "0.57"
$newSheet.Cells.Item(13, 5).Value = "86.30"
$newSheet.Cells.Item(13, 6).Value = "2.96"
$newSheet.Cells.Item(13, 7).Value = "0.0169"
$newSheet.Cells.Item(13, 8).Value = 7
$newSheet.Cells.Item(14, 1).Value = 12
$newSheet.Cells.Item(14, 2).Value = "008262"
$newSheet.Cells.Item(14, 3).Value = "招商研究优选股票C"
$newSheet.Cells.Item(14, 4).Value = "0.54"
$newSheet.Cells.Item(14, 5).Value = "86.08"
$newSheet.Cells.Item(14, 6).Value = "2.57"
$newSheet.Cells.Item(14, 7).Value = "0.0139"
$newSheet.Cells.Item(14, 8).Value = 7
$newSheet.Cells.Item(15, 1).Value = 13
$newSheet.Cells.Item(15, 2).Value = "005091"
$newSheet.Cells.Item(15, 3).Value = "嘉合睿金混合C"
$newSheet.Cells.Item(15, 4).Value = "0.33"
$newSheet.Cells.Item(15, 5).Value = "80.65"
$newSheet.Cells.Item(15, 6).Value = "3.34"
$newSheet.Cells.Item(15, 7).Value = "0.0110"
$newSheet.Cells.Item(15, 8).Value = 9
$newSheet.Cells.Item(16, 1).Value = 14
$newSheet.Cells.Item(16, 2).Value = "011978"
$newSheet.Cells.Item(16, 3).Value = "格林研究优选混合C"
$newSheet.Cells.Item(16, 4).Value = "0.24"
$newSheet.Cells.Item(16, 5).Value = "86.76"
$newSheet.Cells.Item(16, 6).Value = "3.64"
$newSheet.Cells.Item(16, 7).Value = "0.0087"
$newSheet.Cells.Item(16, 8).Value = 10
$newSheet.Cells.Item(17, 1).Value = 15
$newSheet.Cells.Item(17, 2).Value = "005090"
$newSheet.Cells.Item(17, 3).Value = "嘉合睿金混合A"
$newSheet.Cells.Item(17, 4).Value = "0.25"
$newSheet.Cells.Item(17, 5).Value = "80.65"
$newSheet.Cells.Item(17, 6).Value = "3.34"
$newSheet.Cells.Item(17, 7).Value = "0.0084"
$newSheet.Cells.Item(17, 8).Value = 9
$newSheet.Cells.Item(18, 1).Value = 16
$newSheet.Cells.Item(18, 2).Value = "014014"
$newSheet.Cells.Item(18, 3).Value = "招商臻选平衡混合A"
$newSheet.Cells.Item(18, 4).Value = "0.25"
$newSheet.Cells.Item(18, 5).Value = "66.99"
$newSheet.Cells.Item(18, 6).Value = "2.82"
$newSheet.Cells.Item(18, 7).Value = "0.0070"
$newSheet.Cells.Item(18, 8).Value = 7
$newSheet.Cells.Item(19, 1).Value = 17
$newSheet.Cells.Item(19, 2).Value = "003186"
$newSheet.Cells.Item(19, 3).Value = "鹏华兴安定期开放灵活配置混合"
$newSheet.Cells.Item(19, 4).Value = "0.52"
$newSheet.Cells.Item(19, 5).Value = "20.17"
$newSheet.Cells.Item(19, 6).Value = "1.08"
$newSheet.Cells.Item(19, 7).Value = "0.0056"
$newSheet.Cells.Item(19, 8).Value = 6
$newSheet.Cells.Item(20, 1).Value = 18
$newSheet.Cells.Item(20, 2).Value = "014015"
$newSheet.Cells.Item(20, 3).Value = "招商臻选平衡混合C"
$newSheet.Cells.Item(20, 4).Value = "0.19"
$newSheet.Cells.Item(20, 5).Value = "66.99"
$newSheet.Cells.Item(20, 6).Value = "2.82"
$newSheet.Cells.Item(20, 7).Value = "0.0054"
$newSheet.Cells.Item(20, 8).Value = 7
$newSheet.Cells.Item(21, 1).Value = 19
$newSheet.Cells.Item(21, 2).Value = "001664"
$newSheet.Cells.Item(21, 3).Value = "平安鑫安混合A"
$newSheet.Cells.Item(21, 4).Value = "0.33"
$newSheet.Cells.Item(21, 5).Value = "27.32"
$newSheet.Cells.Item(21, 6).Value = "0.88"
$newSheet.Cells.Item(21, 7).Value = "0.0029"
$newSheet.Cells.Item(21, 8).Value = 10
$newSheet.Cells.Item(22, 1).Value = 20
$newSheet.Cells.Item(22, 2).Value = "004791"
$newSheet.Cells.Item(22, 3).Value = "富荣中证500指数增强C"
$newSheet.Cells.Item(22, 4).Value = "0.09"
$newSheet.Cells.Item(22, 5).Value = "90.60"
$newSheet.Cells.Item(22, 6).Value = "2.46"
$newSheet.Cells.Item(22, 7).Value = "0.0022"
$newSheet.Cells.Item(22, 8).Value = 2
$newSheet.Cells.Item(23, 1).Value = 21
$newSheet.Cells.Item(23, 2).Value = "001849"
$newSheet.Cells.Item(23, 3).Value = "前海开源强势共识100强等权重股票"
$newSheet.Cells.Item(23, 4).Value = "0.26"
$newSheet.Cells.Item(23, 5).Value = "71.91"
$newSheet.Cells.Item(23, 6).Value = "0.84"
$newSheet.Cells.Item(23, 7).Value = "0.0022"
$newSheet.Cells.Item(23, 8).Value = 10
$newSheet.Cells.Item(24, 1).Value = 22
$newSheet.Cells.Item(24, 2).Value = "007049"
$newSheet.Cells.Item(24, 3).Value = "平安鑫安混合E"
$newSheet.Cells.Item(24, 4).Value = "0.15"
$newSheet.Cells.Item(24, 5).Value = "27.32"
$newSheet.Cells.Item(24, 6).Value = "0.88"
$newSheet.Cells.Item(24, 7).Value = "0.0013"
$newSheet.Cells.Item(24, 8).Value = 10
$newSheet.Cells.Item(25, 1).Value = 23
$newSheet.Cells.Item(25, 2).Value = "515590"
$newSheet.Cells.Item(25, 3).Value = "前海开源中证500等权重ETF"
$newSheet.Cells.Item(25, 4).Value = "0.34"
$newSheet.Cells.Item(25, 5).Value = "95.19"
$newSheet.Cells.Item(25, 6).Value = "0.31"
$newSheet.Cells.Item(25, 7).Value = "0.0011"
$newSheet.Cells.Item(25, 8).Value = 5
$newSheet.Cells.Item(26, 1).Value = 24
$newSheet.Cells.Item(26, 2).Value = "004790"
$newSheet.Cells.Item(26, 3).Value = "富荣中证500指数增强A"
$newSheet.Cells.Item(26, 4).Value = "0.02"
$newSheet.Cells.Item(26, 5).Value = "90.60"
$newSheet.Cells.Item(26, 6).Value = "2.46"
$newSheet.Cells.Item(26, 7).Value = "0.0005"
$newSheet.Cells.Item(26, 8).Value = 2
$newSheet.Cells.Item(27, 1).Value = 25
$newSheet.Cells.Item(27, 2).Value = "001665"
$newSheet.Cells.Item(27, 3).Value = "平安鑫安混合C"
$newSheet.Cells.Item(27, 4).Value = "0.01"
$newSheet.Cells.Item(27, 5).Value = "27.32"
$newSheet.Cells.Item(27, 6).Value = "0.88"
$newSheet.Cells.Item(27, 7).Value = "0.0001"
$newSheet.Cells.Item(27, 8).Value = 10
$newSheet.Cells.Item(28, 1).Value = 26
$newSheet.Cells.Item(28, 2).Value = "005494"
$newSheet.Cells.Item(28, 3).Value = "鑫元价值精选灵活配置混合C"
$newSheet.Cells.Item(28, 4).Value = "0.00"
$newSheet.Cells.Item(28, 5).Value = "86.30"
$newSheet.Cells.Item(28, 6).Value = "2.96"
$newSheet.Cells.Item(28, 7).Value = 0
$newSheet.Cells.Item(28, 8).Value = 7
